# updated Facebook and Twitter data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B column values (rows 2-34) ---
$bUpdates = @{
    2  = 68908
    3  = 71313
    4  = 76542
    5  = 77505
    6  = 82879
    7  = 73409
    8  = 72919
    9  = 70670
    10 = 68187
    11 = 75491
    12 = 69605
    13 = 54079
    14 = 58837
    15 = 64232
    16 = 77346
    17 = 102216
    18 = 98937
    19 = 80947
    20 = 78113
    21 = 69932
    22 = 70451
    23 = 79590
    24 = 69536
    25 = 64505
    26 = 64146
    27 = 67567
    28 = 75755
    29 = 65241
    30 = 69379
    31 = 69506
    32 = 67714
    33 = 60541
    34 = 62841
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# --- Append new rows 35-42 ---
$newRows = @(
    @{ Row = 35; A = 44500; B = 63794 },
    @{ Row = 36; A = 44530; B = 61712 },
    @{ Row = 37; A = 44561; B = 52151 },
    @{ Row = 38; A = 44592; B = 49053 },
    @{ Row = 39; A = 44620; B = 53953 },
    @{ Row = 40; A = 44651; B = 60374 },
    @{ Row = 41; A = 44681; B = 56871 },
    @{ Row = 42; A = 44712; B = 60566 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = "M"
}
